# Insert a new data row at row 181 (pushing existing rows 181-291 down to
# 182-292) and fill it with the new observation's data. This matches the
# behaviour of the target diff: a weekly price record was added to the
# dataset, the rest of the rows below it simply shift down by one, and the
# sheet dimension grows from A1:R291 to A1:R292.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 181..291 down to 182..292, creating a blank row 181.
$ws.Rows.Item(181).Insert()

# Populate the newly inserted row 181 with the new record.
$ws.Range("A181").Value = 9
$ws.Range("B181").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C181").Value = "Metropolitana"
$ws.Range("D181").Value = 44777
$ws.Range("E181").Value = 13
$ws.Range("F181").Value = 100112001
$ws.Range("G181").Value = "Berenjena"
$ws.Range("H181").Value = "Sin especificar"
$ws.Range("I181").Value = "Primera"
$ws.Range("J181").Value = 90
$ws.Range("K181").Value = 11000
$ws.Range("L181").Value = 12000
$ws.Range("M181").Value = 11500
$ws.Range("N181").Value = "`$/caja 50 unidades"
$ws.Range("O181").Value = "Región de Arica y Parinacota"
$ws.Range("P181").Value = 230
$ws.Range("Q181").Value = 50
$ws.Range("R181").Value = "Hortaliza"

# Make sure column D keeps the date number format used by the rest of the
# column (style index 2 in the original workbook == numFmtId 165).
$ws.Range("D181").NumberFormat = $ws.Range("D182").NumberFormat
